$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsHsr   = $wb.Worksheets.Item("HSR")

# --- About sheet updates ---------------------------------------------------

# Year bumped from 2021 to 2022
$wsAbout.Range("B4").Value = 2022

# Note label: switched from the 2020-Q3/EU28 single-quarter figure to the
# 2019-Q3/4 averaged/EU27 figure
$wsAbout.Range("B7").Value = "2019-Q3/4 savings rate for EU27"

# First line of the explanatory note below the data: "use ... Q3-2020" ->
# "average ... Q3/4-2019" (the following two lines of the note are unchanged)
$wsAbout.Range("A15").Value = "We average data from Q3/4-2019 in part because the household savings"

# Selection moved on the About sheet (no longer the active tab)
$wsAbout.Range("N13").Select()

# --- HSR sheet updates -------------------------------------------------

# Savings rate formula: single hard-coded figure -> average of the two
# (Q3 & Q4 2019) figures
$wsHsr.Range("B2").Formula = "=AVERAGE(0.1197,0.1249)"

# HSR becomes the active sheet/tab, with a new selected cell
$wsHsr.Activate()
$wsHsr.Range("D7").Select()
